# Update receptor / edge weight-specificity values with newly computed TPM-based figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.507281333333333
$ws.Range("N2").Value = 4.521844
$ws.Range("O2").Value = 0.2777659256678391
$ws.Range("P2").Value = 0.2908185693226504
$ws.Range("Q2").Value = 2.422297568672
$ws.Range("R2").Value = 21.800678118048
$ws.Range("S2").Value = 0.2777659256678391
$ws.Range("T2").Value = 0.2908185693226504

# Row 3 (only O3, P3, S3, T3 change)
$ws.Range("O3").Value = 0.2258629392248442
$ws.Range("P3").Value = 0.2364765825413957
$ws.Range("S3").Value = 0.2258629392248442
$ws.Range("T3").Value = 0.2364765825413957

# Row 4
$ws.Range("M4").Value = 0.9373149999999999
$ws.Range("N4").Value = 2.811945
$ws.Range("O4").Value = 0.1727309712259096
$ws.Range("P4").Value = 0.1808478624901656
$ws.Range("Q4").Value = 1.50632519316
$ws.Range("R4").Value = 13.55692673844
$ws.Range("S4").Value = 0.1727309712259096
$ws.Range("T4").Value = 0.1808478624901656

# Row 5
$ws.Range("M5").Value = 0.730656
$ws.Range("N5").Value = 1.461312
$ws.Range("O5").Value = 0.1346472856105345
$ws.Range("P5").Value = 0.0939830443451877
$ws.Range("Q5").Value = 1.174210953984
$ws.Range("R5").Value = 7.045265723904
$ws.Range("S5").Value = 0.1346472856105345
$ws.Range("T5").Value = 0.0939830443451877

# Row 6
$ws.Range("M6").Value = 1.025559333333333
$ws.Range("N6").Value = 3.076678
$ws.Range("O6").Value = 0.1889928782708727
$ws.Range("P6").Value = 0.1978739413006007
$ws.Range("Q6").Value = 1.648139484464
$ws.Range("R6").Value = 14.833255360176
$ws.Range("S6").Value = 0.1889928782708727
$ws.Range("T6").Value = 0.1978739413006007
